$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$sh = $master.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Text = $tr.Text
Write-Output "done"
